# Edit described by the diff:
#  - B2 (teacher name) changes from 张三 to 周福山
#  - A new row 3 is appended with a second teacher's record
#  - Selection cursor moves to K1
#  - Used range / dimension grows to A1:L3 (handled automatically)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2: teacher name ---
$ws.Range("B2").Value = "周福山"

# --- Append new row 3 (second teacher record) ---

# A3: teacher id, stored as text (number format "@" applied before the
# value so it is kept as a text string instead of being parsed as a number)
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2010110156"

# B3: name, C3: gender - plain strings, default format
$ws.Range("B3").Value = "李敏"
$ws.Range("C3").Value = "女"

# D3: password, stored as a genuine number but displayed with a text
# format (value set first, then the display format applied, so it stays
# numeric instead of being converted to text)
$ws.Range("D3").Value = 123456
$ws.Range("D3").NumberFormat = "@"

# E3: birthday, stored as text (format applied first so the date-looking
# string is not auto-converted into a date serial number)
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1975-10-26"

# F3/G3/H3: education, university, political status - plain strings
$ws.Range("F3").Value = "博士"
$ws.Range("G3").Value = "清华大学"
$ws.Range("H3").Value = "党员"

# I3: joining year, genuine number displayed with text format
$ws.Range("I3").Value = 2011
$ws.Range("I3").NumberFormat = "@"

# J3/K3: joining college, title - plain strings
$ws.Range("J3").Value = "计算机科学学院"
$ws.Range("K3").Value = "教授"

# L3: title-acquired date, stored as text
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "2005-10-11"

# --- Move the selection cursor to K1 ---
$ws.Range("K1").Select() | Out-Null
